$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A handful of the new "Price" strings are plain decimal numbers (e.g. "0.999",
# "17.43", "8.43" ...). Assigning those through .Value would make Excel COM
# auto-convert the literal text into a genuine numeric value, losing the exact
# text representation the source data uses (trailing zeros, etc). Forcing the
# cell to Text format ("@") before the assignment keeps it literal; resetting
# the style back to "Normal" afterwards keeps the cell format identical to how
# it started (General) while the stored value remains text.
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"

$ws.Range("D2").Value = "68.225.01"
$ws.Range("E2").Value = "  +0.95%  "
$ws.Range("D3").Value = "3.342.65"
$ws.Range("E3").Value = "  +0.68%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "583.90"
$ws.Range("E5").Value = "  +0.87%  "
$ws.Range("D6").Value = "177.41"
$ws.Range("E6").Value = "  +2.14%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.71%  "
$ws.Range("D9").Value = "0.183"
$ws.Range("E9").Value = "  +4.56%  "
$ws.Range("E10").Value = "  +1.37%  "
$ws.Range("D11").Value = "48.12"
$ws.Range("E11").Value = "  +6.22%  "
$ws.Range("D12").Value = "0.0000273"
$ws.Range("E12").Value = "  +1.89%  "
$ws.Range("D13").Value = "694.33"
$ws.Range("E13").Value = "  +3.98%  "
$ws.Range("D14").Value = "3.883.13"
$ws.Range("D15").Value = "8.43"
$ws.Range("E15").Value = "  +0.71%  "
$ws.Range("D16").Value = "68.308.02"
$ws.Range("E16").Value = "  +0.92%  "
$ws.Range("E17").Value = "  +1.29%  "
$ws.Range("D18").Value = "3.332.05"
$ws.Range("E18").Value = "  +0.52%  "
$ws.Range("D19").Value = "17.43"
$ws.Range("E19").Value = "  +0.23%  "
$ws.Range("E20").Value = "  +2.69%  "
$ws.Range("E21").Value = "  +1.09%  "
$ws.Range("D22").Value = "5.44"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("D23").Value = "16.99"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "100.18"
$ws.Range("E24").Value = "  +2.06%  "
$ws.Range("D25").Value = "3.92"
$ws.Range("E25").Value = "  +2.67%  "
$ws.Range("D26").Value = "2.70"
$ws.Range("E26").Value = "  +1.28%  "
$ws.Range("E27").Value = "  +2.67%  "
$ws.Range("D28").Value = "32.99"
$ws.Range("E28").Value = "  -1.77%  "
$ws.Range("D29").Value = "8.51"
$ws.Range("E29").Value = "  +1.46%  "
$ws.Range("D30").Value = "6.96"
$ws.Range("E30").Value = "  -5.38%  "
$ws.Range("D31").Value = "564.62"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").Value = "11.03"
$ws.Range("E32").Value = "  +1.10%  "
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "57.50"
$ws.Range("E34").Value = "  +1.48%  "
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "3.692.78"
$ws.Range("E36").Value = "  -0.15%  "
$ws.Range("D37").Value = "3.28"
$ws.Range("E37").Value = "  +1.69%  "
$ws.Range("D38").Value = "0.138"
$ws.Range("E38").Value = "  +4.92%  "
$ws.Range("D39").Value = "34.77"
$ws.Range("E39").Value = "  +5.91%  "
$ws.Range("D40").Value = "3.18"
$ws.Range("E40").Value = "  +2.80%  "
$ws.Range("D41").Value = "2.62"
$ws.Range("E41").Value = "  -0.11%  "
$ws.Range("B42").Value = "PEPE"
$ws.Range("C42").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D42").Value = "0.0₃0672"
$ws.Range("E42").Value = "  +2.12%  "
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").Value = "0.335"
$ws.Range("E43").Value = "  +1.01%  "
$ws.Range("D44").Value = "3.29"
$ws.Range("E44").Value = "  +0.81%  "
$ws.Range("E45").Value = "  +2.09%  "
$ws.Range("D46").Value = "2.65"
$ws.Range("E46").Value = "  +2.91%  "
$ws.Range("E47").Value = "  +0.96%  "
$ws.Range("E48").Value = "  -0.01%  "
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("E50").Value = "  +3.13%  "
$ws.Range("E51").Value = "  +0.76%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
